$d = $word.ActiveDocument

$d.Content.Find.Execute("This if efficient", $true, $false, $false, $false, $false,
                         $true, 1, $false, "This is efficient", 2)

